$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.826.27"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "3.377.61"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.374.48"
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "3.943.76"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.38%  "
$ws.Range("E16").Value = "  -4.93%  "
$ws.Range("D17").Value = "3.369.77"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "60.966.60"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").Value = "3.501.09"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("E28").Value = "  +10.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.167"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.771"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "2.359.64"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.78%  "
